$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row of data (row 56) ---

# Date (E56): copy number-format/alignment from E55 (date style), then set
# the date serial value. F column is left untouched here on purpose so the
# SUM(F:F) dependency chain for F56 stays intact.
$ws.Range("E55").Copy()
$ws.Range("E56").PasteSpecial(-4122)
$ws.Range("E56").Value = 43812

# Description (G56): copy formatting (incl. wrap text) from G53, a similar
# multi-line entry, then set the new text.
$ws.Range("G53").Copy()
$ws.Range("G56").PasteSpecial(-4122)
$ws.Range("G56").Value = "Anpassungen der Effekte beim Bewegen über Menüpunkte`nRecherche und Erstimplementierung von ResourceBundles`nCodedokumentierung"

$excel.CutCopyMode = 0

# Hours (F56): plain value assignment - F column's default style already
# matches the desired formatting, and this keeps SUM(F:F) recalculating.
$ws.Range("F56").Value = 5.5

# Row height grows to fit the three lines of wrapped text
$ws.Rows.Item(56).RowHeight = 43.2

# --- Update view state (scroll position / selection) ---
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("G59").Select() | Out-Null

Write-Host "done"
